$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4632500
$ws.Range("C3").Value = 5143192.857142857
$ws.Range("C4").Value = 2835714.285714285
$ws.Range("C5").Value = -3500000
$ws.Range("C6").Value = 3695714.285714286
$ws.Range("C7").Value = 12807121.42857143
